$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues  = -4163
$xlPasteFormats = -4122

# Helper pattern used below for date-like strings (e.g. "2020-08-29"):
# assigning such a string straight to Range.Value gets auto-detected and
# stored as a real date serial. Writing it as a quoted-text formula first
# and then collapsing the formula to its literal value via a
# Copy + PasteSpecial(xlPasteValues) round-trip keeps it as plain text
# (matches the inlineStr text the source file uses) without touching the
# cell's existing style.

# B1: " " -> "akhil"
$ws.Range("B1").Value = "akhil"

# B2: "org name" -> blank
$ws.Range("B2").Value = ""

# B5: Ledger Generation Date "2020-07-30" -> "2020-09-11" (keep as text)
$ws.Range("B5").Formula = '="2020-09-11"'
$ws.Range("B5").Copy()
$ws.Range("B5").PasteSpecial($xlPasteValues)

# Row 10: existing data row updated
$ws.Range("A10").Formula = '="2020-08-29"'
$ws.Range("A10").Copy()
$ws.Range("A10").PasteSpecial($xlPasteValues)
$ws.Range("B10").Value = "test"
$ws.Range("C10").Value = "CN-0001"
$ws.Range("D10").Value = "CREDITNOTE"
$ws.Range("E10").Value = 0

# Row 11: existing data row updated
$ws.Range("A11").Formula = '="2020-09-10"'
$ws.Range("A11").Copy()
$ws.Range("A11").PasteSpecial($xlPasteValues)
$ws.Range("B11").Value = "test"
$ws.Range("C11").Value = "Exp-1"
$ws.Range("D11").Value = "EXPENSE"
$ws.Range("E11").Value = 25

# Insert a new row at 12 (old TOTAL row shifts down to 13); copy the
# formatting used by row 11 so the new data row matches rows 10/11.
$ws.Rows("12").Insert()
$ws.Range("A11:E11").Copy()
$ws.Range("A12:E12").PasteSpecial($xlPasteFormats)

$ws.Range("A12").Formula = '="2020-09-10"'
$ws.Range("A12").Copy()
$ws.Range("A12").PasteSpecial($xlPasteValues)
$ws.Range("B12").Value = "test"
$ws.Range("C12").Value = "Exp-2"
$ws.Range("D12").Value = "EXPENSE"
$ws.Range("E12").Value = 5

# Row 13 (formerly row 12): TOTAL row, formula now spans through row 12
$ws.Range("E13").Formula = "=SUM(E10:E12)"
